$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1833
$ws.Range("I28").Value = 1142.2307
$ws.Range("J28").Value = 4826.3335
$ws.Range("K28").Value = 1142.2307
$ws.Range("L28").Value = 4826.3335
$ws.Range("M28").Value = -657.2307000000001
$ws.Range("N28").Value = -5796.3335
$ws.Range("H76").Value = 2975
$ws.Range("I76").Value = 2975
$ws.Range("K76").Value = 2975
$ws.Range("M76").Value = -2660
$ws.Range("H79").Value = 2975
$ws.Range("I79").Value = 2975
$ws.Range("K79").Value = 2975
$ws.Range("M79").Value = -1883
$ws.Range("H137").Value = 1535.2195
$ws.Range("I137").Value = 1145.8928
$ws.Range("J137").Value = 2373.7693
$ws.Range("K137").Value = 3437.6784
$ws.Range("L137").Value = 7121.3079
$ws.Range("M137").Value = -887.6784000000002
$ws.Range("N137").Value = -12221.3079
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 18799.2
$ws.Range("J52").Value = 18799.2
$ws.Range("L52").Value = 18799.2
$ws.Range("N52").Value = -19435.2
$ws.Range("H61").Value = 1047.3055
$ws.Range("I61").Value = 903.43335
$ws.Range("J61").Value = 1766.6666
$ws.Range("K61").Value = 903.43335
$ws.Range("L61").Value = 1766.6666
$ws.Range("M61").Value = -691.43335
$ws.Range("N61").Value = -2190.6666
$ws.Range("H63").Value = 3916.6667
$ws.Range("I63").Value = 2833.3333
$ws.Range("K63").Value = 2833.3333
$ws.Range("M63").Value = -2147.3333
$ws.Range("H66").Value = 3916.6667
$ws.Range("I66").Value = 2833.3333
$ws.Range("K66").Value = 14166.6665
$ws.Range("M66").Value = -10734.6665
$ws.Range("H80").Value = 21275
$ws.Range("I80").Value = 5100
$ws.Range("J80").Value = 26666.666
$ws.Range("K80").Value = 5100
$ws.Range("L80").Value = 26666.666
$ws.Range("M80").Value = -4102
$ws.Range("N80").Value = -28662.666
$ws.Range("H83").Value = 21275
$ws.Range("I83").Value = 5100
$ws.Range("J83").Value = 26666.666
$ws.Range("K83").Value = 15300
$ws.Range("L83").Value = 79999.99800000001
$ws.Range("M83").Value = -10308
$ws.Range("N83").Value = -89983.99800000001
$ws.Range("H132").Value = 2119.8286
$ws.Range("I132").Value = 2079.24
$ws.Range("J132").Value = 2221.3
$ws.Range("K132").Value = 6237.719999999999
$ws.Range("L132").Value = 6663.900000000001
$ws.Range("M132").Value = -3707.719999999999
$ws.Range("N132").Value = -11723.9
$ws.Range("H136").Value = 1047.3055
$ws.Range("I136").Value = 903.43335
$ws.Range("J136").Value = 1766.6666
$ws.Range("K136").Value = 2710.30005
$ws.Range("L136").Value = 5299.9998
$ws.Range("M136").Value = -160.3000499999998
$ws.Range("N136").Value = -10399.9998
$ws.Range("H137").Value = 39600
$ws.Range("J137").Value = 39600
$ws.Range("L137").Value = 39600
$ws.Range("N137").Value = -49800
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29199.666
$ws.Range("J51").Value = 29199.666
$ws.Range("L51").Value = 29199.666
$ws.Range("N51").Value = -30181.666
$ws.Range("H57").Value = 43633.332
$ws.Range("J57").Value = 43633.332
$ws.Range("L57").Value = 43633.332
$ws.Range("N57").Value = -45073.332
$ws.Range("H134").Value = 566628.4399999999
$ws.Range("I134").Value = 871931.75
$ws.Range("J134").Value = 4870.36
$ws.Range("K134").Value = 2615795.25
$ws.Range("L134").Value = 14611.08
$ws.Range("M134").Value = -2613260.25
$ws.Range("N134").Value = -19681.08
$ws.Range("H136").Value = 43633.332
$ws.Range("J136").Value = 43633.332
$ws.Range("L136").Value = 43633.332
$ws.Range("N136").Value = -53833.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14085769
$ws.Range("I31").Value = 19231532
$ws.Range("J31").Value = 2627.2104
$ws.Range("K31").Value = 19231532
$ws.Range("L31").Value = 2627.2104
$ws.Range("M31").Value = -19231237
$ws.Range("N31").Value = -3217.2104
$ws.Range("H34").Value = 14085769
$ws.Range("I34").Value = 19231532
$ws.Range("J34").Value = 2627.2104
$ws.Range("K34").Value = 19231532
$ws.Range("L34").Value = 2627.2104
$ws.Range("M34").Value = -19231330
$ws.Range("N34").Value = -3031.2104
$ws.Range("H58").Value = 2719.3215
$ws.Range("I58").Value = 2887.7346
$ws.Range("J58").Value = 1540.4286
$ws.Range("K58").Value = 2887.7346
$ws.Range("L58").Value = 1540.4286
$ws.Range("M58").Value = -2684.7346
$ws.Range("N58").Value = -1946.4286
$ws.Range("H134").Value = 17188980
$ws.Range("I134").Value = 1819566.9
$ws.Range("K134").Value = 5458700.699999999
$ws.Range("M134").Value = -5456165.699999999
$ws.Range("H136").Value = 2719.3215
$ws.Range("I136").Value = 2887.7346
$ws.Range("J136").Value = 1540.4286
$ws.Range("K136").Value = 8663.203799999999
$ws.Range("L136").Value = 4621.2858
$ws.Range("M136").Value = -6113.203799999999
$ws.Range("N136").Value = -9721.2858
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13889661
$ws.Range("I5").Value = 30303488
$ws.Range("J5").Value = 1038.7693
$ws.Range("K5").Value = 90910464
$ws.Range("L5").Value = 3116.3079
$ws.Range("M5").Value = -90910352
$ws.Range("N5").Value = -3340.3079
$ws.Range("H80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4872
$ws.Range("H83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("L83").Value = 9000
$ws.Range("N83").Value = -18360
$ws.Range("H122").Value = 586.5333000000001
$ws.Range("I122").Value = 518.9545000000001
$ws.Range("J122").Value = 772.375
$ws.Range("K122").Value = 4670.5905
$ws.Range("L122").Value = 6951.375
$ws.Range("M122").Value = -2220.5905
$ws.Range("N122").Value = -11851.375
$ws.Range("H131").Value = 914.34
$ws.Range("I131").Value = 875
$ws.Range("J131").Value = 915.9792
$ws.Range("K131").Value = 2625
$ws.Range("L131").Value = 2747.9376
$ws.Range("M131").Value = 2415
$ws.Range("N131").Value = -12827.9376
$ws.Range("H135").Value = 13889661
$ws.Range("I135").Value = 30303488
$ws.Range("J135").Value = 1038.7693
$ws.Range("K135").Value = 272731392
$ws.Range("L135").Value = 9348.923699999999
$ws.Range("M135").Value = -272728857
$ws.Range("N135").Value = -14418.9237
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 134.11765
$ws.Range("I2").Value = 31.6
$ws.Range("J2").Value = 176.83333
$ws.Range("K2").Value = 31.6
$ws.Range("L2").Value = 176.83333
$ws.Range("M2").Value = 81.40000000000001
$ws.Range("N2").Value = -402.83333
$ws.Range("H57").Value = 14570.353
$ws.Range("J57").Value = 15387.25
$ws.Range("L57").Value = 15387.25
$ws.Range("N57").Value = -17027.25
$ws.Range("H62").Value = 29750
$ws.Range("J62").Value = 29750
$ws.Range("L62").Value = 29750
$ws.Range("N62").Value = -31122
$ws.Range("H65").Value = 29750
$ws.Range("J65").Value = 29750
$ws.Range("L65").Value = 89250
$ws.Range("N65").Value = -96114
$ws.Range("H132").Value = 2634201
$ws.Range("I132").Value = 2879.05
$ws.Range("J132").Value = 5557892
$ws.Range("K132").Value = 8637.150000000001
$ws.Range("L132").Value = 16673676
$ws.Range("M132").Value = -6107.150000000001
$ws.Range("N132").Value = -16678736
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 555.8421
$ws.Range("I55").Value = 136.88889
$ws.Range("J55").Value = 932.9
$ws.Range("K55").Value = 136.88889
$ws.Range("L55").Value = 932.9
$ws.Range("M55").Value = 36.11111
$ws.Range("N55").Value = -1278.9
$ws.Range("H64").Value = 32800
$ws.Range("J64").Value = 32800
$ws.Range("L64").Value = 32800
$ws.Range("N64").Value = -33250
$ws.Range("H67").Value = 32800
$ws.Range("J67").Value = 32800
$ws.Range("L67").Value = 32800
$ws.Range("N67").Value = -34360
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -98112
$ws.Range("H132").Value = 4247.613
$ws.Range("I132").Value = 5583.778
$ws.Range("J132").Value = 2397.5386
$ws.Range("K132").Value = 16751.334
$ws.Range("L132").Value = 7192.6158
$ws.Range("M132").Value = -14221.334
$ws.Range("N132").Value = -12252.6158
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2500.451
$ws.Range("I132").Value = 2730.7144
$ws.Range("J132").Value = 1996.75
$ws.Range("K132").Value = 8192.143199999999
$ws.Range("L132").Value = 5990.25
$ws.Range("M132").Value = -5662.143199999999
$ws.Range("N132").Value = -11050.25
